$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we are about to write keep their original "text" nature
# (the source workbook stores these values as inline strings, not numbers),
# so force the Text number format before assigning values.
$ws.Range("B2:E51").NumberFormat = "@"

$updates = @{
    2  = @{ D = "245.73"; E = "-0.43%" }
    3  = @{ D = "30.11"; E = "-0.40%" }
    4  = @{ D = "5.155"; E = "-0.39%" }
    5  = @{ D = "0.05759" }
    6  = @{ D = "6.669"; E = "0.98%" }
    7  = @{ D = "3.283"; E = "6.95%" }
    8  = @{ D = "0.8495"; E = "-0.83%" }
    9  = @{ D = "0.8584"; E = "-2.44%" }
    10 = @{ D = "0.1383"; E = "1.15%" }
    11 = @{ B = "MandalaExchangeToken"; C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D = "0.07083"; E = "0.12%" }
    12 = @{ B = "BitrueCoin"; C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D = "0.03239"; E = "12.95%" }
    13 = @{ B = "BitMartToken"; C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D = "0.09362"; E = "-0.31%" }
    14 = @{ B = "BitForexToken"; C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D = "0.001535"; E = "0.80%" }
    15 = @{ B = "TigerCash"; C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D = "0.005921"; E = "-2.44%" }
    16 = @{ B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "3.547"; E = "1.62%" }
    17 = @{ B = "BTSEToken"; C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D = "2.218"; E = "-2.58%" }
    18 = @{ B = "BitpandaEcosystemToken"; C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D = "0.3123"; E = "-1.93%" }
    19 = @{ B = "LiechtensteinCryptoassetsExchange"; C = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D = "0.03422"; E = "5.24%" }
    20 = @{ B = "ProBitToken"; C = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; D = "0.1316"; E = "1.17%" }
    21 = @{ B = "MCDex"; C = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; D = "3.483"; E = "-0.81%" }
    22 = @{ B = "ZBToken"; C = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; D = "0.1410"; E = "2.15%" }
    23 = @{ B = "CoinExToken"; C = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D = "0.04125"; E = "-0.42%" }
    24 = @{ B = "One"; C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D = "0.0005960"; E = "-94.22%" }
    25 = @{ E = "1.02%" }
    26 = @{ E = "-7.55%" }
    27 = @{ E = "-0.82%" }
    40 = @{ D = "0.03755"; E = "-0.77%" }
    41 = @{ D = "0.1069"; E = "-0.06%" }
    42 = @{ D = "0.002200"; E = "0.01%" }
    43 = @{ D = "0.002950"; E = "-48.51%" }
    44 = @{ D = "0.009454"; E = "-5.97%" }
    45 = @{ D = "0.00005497"; E = "8.01%" }
    46 = @{ E = "0.01%" }
    47 = @{ D = "0.07100" }
    49 = @{ E = "0.01%" }
    50 = @{ E = "0.01%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
